$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New cell Sheet1!B4 holding a UTF-8 (Hindi) greeting, styled like the rest of
# the table (center-aligned, matching cellXf index 1 used by A1:C3).
$ws.Range("B4").Value = "नमस्ते"
$ws.Range("B4").HorizontalAlignment = -4108

# Make B4 the active selection on Sheet1 (matches the recorded <selection/>).
$ws.Range("B4").Select()

# New workbook-scoped defined name using a non-ASCII (Greek) identifier that
# points at the new cell. The COM name-creation call only accepts ASCII
# literal names, so create it under a throwaway ASCII name first and then
# rename the Name object to the Unicode target - renaming isn't subject to
# the same restriction.
$wb.Names.Add('TempGreekName', '=Sheet1!$B$4')
$tempName = $wb.Names.Item('TempGreekName')
$tempName.Name = 'Χαιρετισμός'
